# Rename the worksheet from "water_use" to "water_use_statistics"
# (commit message: "Updated to point at mapserver instead of fs")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "water_use_statistics"

# Move the active window's on-screen position (matches the xWindow/yWindow
# change recorded for the workbook view).
$aw = $excel.ActiveWindow
$aw.Left = 28680
$aw.Top = 525

# Update the selected/active cell on the sheet from M22 to G17.
$ws.Range("G17").Select()
